$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650996150159079"
$ws1.Range("B2").Value = "go_stims-16509961501270394.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961501430407.csv"
$ws1.Range("B4").Value = "go_stims-16509961501430407.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996150159079.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509961533910403"
$ws2.Range("B2").Value = "OB-16509961504310482.csv"
$ws2.Range("B3").Value = "ZB-match_1-1650996150311078.csv"
$ws2.Range("B4").Value = "OB-16509961506550412.csv"
$ws2.Range("B5").Value = "TB-1650996151663043.csv"
$ws2.Range("B6").Value = "ZB-match_6-16509961501750405.csv"
$ws2.Range("B7").Value = "TB-16509961533590755.csv"
$ws2.Range("B8").Value = "TB-16509961520950756.csv"
$ws2.Range("B9").Value = "OB-16509961503350809.csv"
$ws2.Range("B10").Value = "ZB-match_4-16509961501990833.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509961533910403"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961534390776"
$ws4.Range("B2").Value = "MM_stims-16509961534070473.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961533910403.csv"
$ws4.Range("B4").Value = "MM_stims-16509961534230664.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961534070473.csv"
$ws4.Range("B6").Value = "MM_stims-16509961534390776.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961534230664.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961535030794"
$ws5.Range("B2").Value = "SAT_stims-16509961534550803.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961534390776.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961534870791.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961534710886.csv"
